$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural shape fixes first (row count + which rows carry B/C) ---

# The table shrinks from 22 to 21 data rows; drop the old trailing row.
$ws.Rows.Item(22).Delete()

# Row 17 goes from a full A/B/C row to an A-only row (like rows 12 and 18
# used to be). Deleting it shifts rows 18+ up by one, which turns the old
# row 18 (A-only) into the new row 17 for free; re-inserting a blank row
# at 18 keeps every later row number lined up with the target layout.
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(18).Insert()

# --- Now (re)write every cells final text, matching the target table ---

# Row 10
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "849935 - Humberto Felipe da Silva"
$ws.Range("C10").Value = "849935 - Humberto Felipe da Silva"
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Lead students to understand the key innovation management issues, their characteristics and critical points for success, through real business problems. The activities will be carried out as a team and will focus on the development of the necessary skills to successfully manage the entire innovation process, from its conception to placing the product on the market"
$ws.Range("C11").Value = "Lead students to understand the key innovation management issues, their characteristics and critical points for success, through real business problems. The activities will be carried out as a team and will focus on the development of the necessary skills to successfully manage the entire innovation process, from its conception to placing the product on the market"
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Rows.Item(12).AutoFit()

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Innovation management. Costs of innovation. Innovation implementation process. Transforming ideation into business. The Innovation Market"
$ws.Range("C14").Value = "Innovation management. Costs of innovation. Innovation implementation process. Transforming ideation into business. The Innovation Market"
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2020"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2020"
$ws.Rows.Item(15).RowHeight = 120
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Which leads some organizations to have the necessity to manage the innovation; what are the main factors driving innovation and how does the diffusion process occur. What are the main types of innovation that need to be considered; What are the main strategies to implement one of innovation in the market; what and how companies seek innovative ideas; What and how different factors influence the way managers to prioritize innovation choices; What are the key challenges in putting innovations into practice? How managers can build an organization focused on innovation as a market strategy; what are the main success evaluation systems of an innovative idea; Resource management and innovative programs in a company."
$ws.Range("C16").Value = "Which leads some organizations to have the necessity to manage the innovation; what are the main factors driving innovation and how does the diffusion process occur. What are the main types of innovation that need to be considered; What are the main strategies to implement one of innovation in the market; what and how companies seek innovative ideas; What and how different factors influence the way managers to prioritize innovation choices; What are the key challenges in putting innovations into practice? How managers can build an organization focused on innovation as a market strategy; what are the main success evaluation systems of an innovative idea; Resource management and innovative programs in a company."
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "849935 - Humberto Felipe da Silva"
$ws.Range("C18").Value = "849935 - Humberto Felipe da Silva"
$ws.Rows.Item(18).RowHeight = 60
$ws.Range("B10").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Seminários e Estudos de Casos, aplicação de “Pitchs” (breve apresentação oral de uma ideia, produto ou oportunidade de negócio) e outras formas de apresentação de ideias em empresas"
$ws.Range("C19").Value = "Seminários e Estudos de Casos, aplicação de “Pitchs” (breve apresentação oral de uma ideia, produto ou oportunidade de negócio) e outras formas de apresentação de ideias em empresas"
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "As avaliações serão: a) contínuas considerando a participação dos alunos nas atividades; b) avaliação das apresentações parciais dos trabalhos; e c) apresentação final dos trabalhos."
$ws.Range("C20").Value = "As avaliações serão: a) contínuas considerando a participação dos alunos nas atividades; b) avaliação das apresentações parciais dos trabalhos; e c) apresentação final dos trabalhos."
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina."
$ws.Range("C21").Value = "Reapresentação do último seminário, cuja nota constituirá a nota final da disciplina."
$ws.Rows.Item(21).RowHeight = 120
